$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the affected cells are written back as text (matching the
# original t="str" cell type) rather than being auto-converted to
# numbers by Excel's value-type inference.
$targets = @(
    @{ Cell = "C2"; Value = "0" },
    @{ Cell = "D2"; Value = "2" },

    @{ Cell = "C3"; Value = "2" },
    @{ Cell = "D3"; Value = "4" },

    @{ Cell = "C4"; Value = "1" },
    @{ Cell = "D4"; Value = "2" },

    @{ Cell = "C5"; Value = "0" },
    @{ Cell = "D5"; Value = "1" },
    @{ Cell = "E5"; Value = "0" },

    @{ Cell = "C6"; Value = "9" },
    @{ Cell = "D6"; Value = "6" },
    @{ Cell = "E6"; Value = "1" }
)

foreach ($t in $targets) {
    $rng = $ws.Range($t.Cell)
    $rng.NumberFormat = "@"
    $rng.Value = $t.Value
}
